$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert 10 new rows starting at row 27, pushing existing rows 27+ down by 10
$ws.Range("A27:A36").EntireRow.Insert()

# Fill the newly inserted rows 27-29 with problem 11 (common part) data.
# Columns A and B are filled first (row by row), then column C is filled in a
# second pass -- this mirrors the shared-string insertion order in the target
# workbook (codes/sentences interned before the keyword column).
$ws.Range("A27").Value = "c0017"
$ws.Range("B27").Value = "주어진 삼각함수의 그래프 개형을 파악하기 위해서 삼각함수의 주기를 구합니다."

$ws.Range("A28").Value = "c0018"
$ws.Range("B28").Value = "두 꼭짓점의 좌표로 부터 정삼각형의 한 변의 길이를 알아냅니다."

$ws.Range("A29").Value = "c0019"
$ws.Range("B29").Value = "두 변의 길이와 끼인 각의 사인값을 이용해서 삼각형의 넓이를 구합니다."

$ws.Range("C27").Value = "탄젠트함수;"
$ws.Range("C28").Value = "두 점 사이의 거리;"
$ws.Range("C29").Value = "정삼각형의 넓이;"

# Column widths (engine snaps ColumnWidth to 1/7-character increments, same as
# Excel's internal pixel-based storage grid; 98.5714... / 64.2857... are the
# inputs whose stored widths land closest to / exactly on the target widths
# of 99.25 and 65 characters respectively)
$ws.Columns.Item(2).ColumnWidth = 98.57142857142857
$ws.Columns.Item(3).ColumnWidth = 64.28571428571429

# Update selection to match the final edited cell
$ws.Range("C29").Select()
